# Weekly update: a new price observation for "Vega Modelo de Temuco -
# Frambuesa" is inserted as the new row 3 (week of 44592), pushing the
# previously-existing rows 3..9 down to rows 4..10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 3 (shifts old rows 3-9 down to 4-10,
# inheriting formatting/styles from the row above, same as Excel's UI).
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with this week's observation.
$ws.Cells.Item(3, 1).Value  = 10
$ws.Cells.Item(3, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(3, 3).Value  = "La Araucanía"
$ws.Cells.Item(3, 4).Value  = 44592
$ws.Cells.Item(3, 5).Value  = 9
$ws.Cells.Item(3, 6).Value  = "Fruta"
$ws.Cells.Item(3, 7).Value  = 100101
$ws.Cells.Item(3, 8).Value  = "Berries"
$ws.Cells.Item(3, 9).Value  = 100101004
$ws.Cells.Item(3, 10).Value = "Frambuesa"
$ws.Cells.Item(3, 11).Value = "Sin especificar"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 5
$ws.Cells.Item(3, 14).Value = 7500
$ws.Cells.Item(3, 15).Value = 7500
$ws.Cells.Item(3, 16).Value = 7500
$ws.Cells.Item(3, 17).Value = "$/envase 1 kilo"
$ws.Cells.Item(3, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(3, 19).Value = 7500
$ws.Cells.Item(3, 20).Value = 1
